# Word document currently has (within the paragraph containing the unique
# word "principallem"):
#
#   ... principallem<exp>ent</exp><lb/> ...
#
# Target change:
#   1. Insert a new red "<corr>" run (Courier New, 9pt) immediately before
#      the existing "<exp>" run.
#   2. Drop the explicit black color on the "ent" run (falls back to
#      automatic/inherited color).
#   3. Insert a new red "</corr>" run (Courier New, 9pt) immediately after
#      the existing "</exp>" run.

$d = $word.ActiveDocument

# "principallem" occurs exactly once in the document, so anchor on it to
# reach the right paragraph (the text "<exp>ent</exp>" also appears earlier
# in the document, attached to a different word, so a plain search for
# "<exp>" alone would be ambiguous).
$anchor = $d.Content
$anchor.Find.Execute("principallem", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$searchStart = $anchor.End

# Locate the "<exp>" run that follows.
$expTag = $d.Range($searchStart, $d.Content.End)
$expTag.Find.Execute("<exp>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# --- Step 1: insert "<corr>" right before "<exp>" -------------------------
# Copy the "<lb/>" run later in the same paragraph as a formatting donor: it
# already carries Courier New / sz 18 / szCs 18, matching what "<corr>"
# needs -- only its color (and text) must change afterwards. Pasting a
# duplicate run (rather than building one from scratch through Font.*)
# keeps the full rFonts/szCs detail that a brand new run would otherwise
# lose.
$lbDonor = $d.Range($searchStart, $d.Content.End)
$lbDonor.Find.Execute("<lb/>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$lbDonor.Copy()

$corrInsertPos = $expTag.Start
$corrPaste = $d.Range($corrInsertPos, $corrInsertPos)
$corrPaste.Paste()

$corrRun = $d.Range($corrInsertPos, $corrInsertPos + 5)   # length of "<lb/>"
$corrRun.Text = "<corr>"
$corrRun = $d.Range($corrInsertPos, $corrInsertPos + 6)   # length of "<corr>"
$corrRun.Font.Color = 1118633                              # RGB a91111 -> BGR int

# --- Step 2: clear the explicit black color on "ent" -----------------------
$expTag2 = $d.Range($corrRun.End, $d.Content.End)
$expTag2.Find.Execute("<exp>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$entRun = $d.Range($expTag2.End, $expTag2.End + 3)          # "ent"
$entRun.Font.Color = -16777216                              # wdColorAutomatic

# --- Step 3: insert "</corr>" right after "</exp>" -------------------------
$closeExpTag = $d.Range($entRun.End, $d.Content.End)
$closeExpTag.Find.Execute("</exp>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$lbDonor2 = $d.Range($closeExpTag.End, $d.Content.End)
$lbDonor2.Find.Execute("<lb/>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$lbDonor2.Copy()

$closeCorrInsertPos = $closeExpTag.End
$closeCorrPaste = $d.Range($closeCorrInsertPos, $closeCorrInsertPos)
$closeCorrPaste.Paste()

$closeCorrRun = $d.Range($closeCorrInsertPos, $closeCorrInsertPos + 5)  # "<lb/>"
$closeCorrRun.Text = "</corr>"
$closeCorrRun = $d.Range($closeCorrInsertPos, $closeCorrInsertPos + 7)  # "</corr>"
$closeCorrRun.Font.Color = 1118633

Write-Output "Edit applied"
